# ADDED NA'S TO DATAFILE
# Insert a new "county" variable row into the metaware_meta_clean.csv codebook sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metaware_meta_clean.csv")
$ws.Activate()

# Insert a new row at position 9 (pushes existing rows 9-28 down to 10-29),
# inheriting the formatting of the row above it.
$ws.Rows("9:9").Insert() | Out-Null

# Populate the new row with the "county" codebook entry.
$ws.Range("A9").Value = "county"
$ws.Range("B9").Value = "Indicator of what country the investigation was conducted in (or, if not stated, the country of the corresponding author)"
$ws.Range("C9").Value = "coded according to ISO-3 standards"
$ws.Rows("9:9").RowHeight = 16

# Update the view's selection to match where the editor left off.
$ws.Range("A17").Select() | Out-Null

Write-Host "Inserted 'county' codebook row into metaware_meta_clean.csv"
